$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New roster data: Player, Position, Team
$data = @(
    @("Jalen Brunson", "PG", "New York Knicks"),
    @("Coby White", "PG,SG", "Chicago Bulls"),
    @("Kelly Oubre Jr.", "SG,SF", "Philadelphia 76ers"),
    @("Caris LeVert", "SG,SF", "Cleveland Cavaliers"),
    @("Norman Powell", "SG,SF", "LA Clippers"),
    @("Jeremy Sochan", "SF,PF", "San Antonio Spurs"),
    @("Devin Vassell", "SG,SF", "San Antonio Spurs"),
    @("Desmond Bane", "SG,SF", "Memphis Grizzlies"),
    @("Yves Missi", "C", "New Orleans Pelicans"),
    @("Walker Kessler", "C", "Utah Jazz"),
    @("LeBron James", "SF,PF", "Los Angeles Lakers"),
    @("Trae Young", "PG", "Atlanta Hawks"),
    @("Guerschon Yabusele", "PF,C", "Philadelphia 76ers"),
    @("Alperen Sengün", "C", "Houston Rockets"),
    @("Jabari Smith Jr.", "PF,C", "Houston Rockets"),
    @("Immanuel Quickley", "PG,SG", "Toronto Raptors"),
    @("Kawhi Leonard", "SG,SF,PF", "LA Clippers"),
    @("Devin Booker", "PG,SG", "Phoenix Suns")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
